# Daily attendance processing - swap the order of "System" and the
# recorder's email address in the "Recorded By" column (column G).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$first = $ws.Cells.Find($oldValue)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $addresses = @()
    $current = $first
    do {
        $addresses += $current.Address()
        $current = $ws.Cells.FindNext($current)
    } while ($current -ne $null -and $current.Address() -ne $firstAddress)

    foreach ($address in $addresses) {
        $ws.Range($address).Value2 = $newValue
    }
}
